# Sprint0/PG3_S0.docx — "Added numpy and fixed outputs"
#
# Semantic changes applied (see commit message / xml diff):
#   1. The "Can be changed to complete user input (string)" bullet gains
#      "for month, day, year " before "(string)".
#   2. The final (previously empty) bullet under "Major abstractions and
#      relationships" is filled with the "holidays" sentence, and a new
#      bullet is appended after it with the "percentages" sentence.
#   3. The hidden "_GoBack" bookmark — which Word re-stamps at the most
#      recently edited spot — moves from the old "numpy, pandas" run to
#      the very end of the document (after the new last bullet).

$d = $word.ActiveDocument

# --- 1. "Can be changed to complete user input (string)" -> "... for month, day, year (string)" ---
$findRange = $d.Content
$found = $findRange.Find.Execute("user input ", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if ($found) {
    $findRange.Collapse(0)
    $findRange.InsertAfter("for month, day, year ")
}

# --- 2. Drop the _GoBack bookmark from its old spot (end of "numpy, pandas" run) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 3. Turn the trailing empty bullet into the "holidays" sentence ---
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "The holidays used in this software will be pulled from USA holidays from holidays package in python"

# --- 4. Append a new sibling bullet (same list style/level) with the "percentages" sentence ---
$lastPara.Range.InsertParagraphAfter()
$newLastPara = $d.Paragraphs.Last
# Type a sentinel trailing character so the bookmark insertion point below isn't the
# paragraph's very-last-character position (collapsed ranges landing exactly before a
# paragraph mark get mis-resolved), then strip the sentinel back out afterwards.
$newLastPara.Range.Text = "Certain percentages will be guessed and used. These will be provided if wanted and can be changed easilyX"

# --- 5. Re-create _GoBack at the new end of the document, right after the text ---
$anchorPos = $newLastPara.Range.End - 2
$anchorRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $anchorRange)

$sentinelRange = $d.Range($newLastPara.Range.End - 2, $newLastPara.Range.End - 1)
$sentinelRange.Delete()

Write-Output "Applied: numpy/outputs fix + holidays/percentages bullets + _GoBack relocation"
